# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text block (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText -replace [regex]::Escape("✅ 1000 Bs = 12.55 = 50288.21 pesos"), "✅ 1000 Bs = 12.62 = 50513.06 pesos"
$newText = $newText -replace [regex]::Escape("✅ 50288.21 pesos = 12.52 = 970.46 Bs"), "✅ 50513.06 pesos = 12.58 = 981.32 Bs"
$cellA1.Value2 = $newText

# --- tasas: update the rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 79.26900000000001
$wsTasas.Range("O10").Value = 4004.12
$wsTasas.Range("N12").Value = 4015
$wsTasas.Range("O12").Value = 78
